$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.567946
$ws.Range("H2").Value = 7.703838
$ws.Range("I2").Value = 0.8361514603436194
$ws.Range("J2").Value = 0.8361514603436195
$ws.Range("M2").Value = 0.4652636666666667
$ws.Range("N2").Value = 1.395791
$ws.Range("O2").Value = 0.02604271297411062
$ws.Range("P2").Value = 0.02604271297411062
$ws.Range("Q2").Value = 1.194771971762
$ws.Range("R2").Value = 10.752947745858
$ws.Range("S2").Value = 0.02177565248461232
$ws.Range("T2").Value = 0.02177565248461232

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.567946
$ws.Range("H3").Value = 7.703838
$ws.Range("I3").Value = 0.8361514603436194
$ws.Range("J3").Value = 0.8361514603436195
$ws.Range("O3").Value = 0.09971126509087273
$ws.Range("P3").Value = 0.09971126509087272
$ws.Range("Q3").Value = 4.574493637353999
$ws.Range("R3").Value = 41.170442736186
$ws.Range("S3").Value = 0.083373719918443
$ws.Range("T3").Value = 0.083373719918443

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.567946
$ws.Range("H4").Value = 7.703838
$ws.Range("I4").Value = 0.8361514603436194
$ws.Range("J4").Value = 0.8361514603436195
$ws.Range("M4").Value = 15.618761
$ws.Range("N4").Value = 46.856283
$ws.Range("O4").Value = 0.8742460219350168
$ws.Range("P4").Value = 0.8742460219350167
$ws.Range("Q4").Value = 40.108134834906
$ws.Range("R4").Value = 360.973213514154
$ws.Range("S4").Value = 0.7310020879405642
$ws.Range("T4").Value = 0.7310020879405642

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5032033333333333
$ws.Range("H5").Value = 1.50961
$ws.Range("I5").Value = 0.1638485396563805
$ws.Range("J5").Value = 0.1638485396563805
$ws.Range("M5").Value = 0.4652636666666667
$ws.Range("N5").Value = 1.395791
$ws.Range("O5").Value = 0.02604271297411062
$ws.Range("P5").Value = 0.02604271297411062
$ws.Range("Q5").Value = 0.2341222279455556
$ws.Range("R5").Value = 2.10710005151
$ws.Range("S5").Value = 0.0042670604894983
$ws.Range("T5").Value = 0.0042670604894983

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5032033333333333
$ws.Range("H6").Value = 1.50961
$ws.Range("I6").Value = 0.1638485396563805
$ws.Range("J6").Value = 0.1638485396563805
$ws.Range("O6").Value = 0.09971126509087273
$ws.Range("P6").Value = 0.09971126509087272
$ws.Range("Q6").Value = 0.8963975280744444
$ws.Range("R6").Value = 8.067577752669999
$ws.Range("S6").Value = 0.01633754517242973
$ws.Range("T6").Value = 0.01633754517242973

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5032033333333333
$ws.Range("H7").Value = 1.50961
$ws.Range("I7").Value = 0.1638485396563805
$ws.Range("J7").Value = 0.1638485396563805
$ws.Range("M7").Value = 15.618761
$ws.Range("N7").Value = 46.856283
$ws.Range("O7").Value = 0.8742460219350168
$ws.Range("P7").Value = 0.8742460219350167
$ws.Range("Q7").Value = 7.859412597736666
$ws.Range("R7").Value = 70.73471337963001
$ws.Range("S7").Value = 0.1432439339944525
$ws.Range("T7").Value = 0.1432439339944525
